$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 21:46"

# Estados Unidos (row 6) - updated case counts
$ws.Range("B6").Value = 52921
$ws.Range("C6").Value = 9187
$ws.Range("E6").Value = 51867
$ws.Range("G6").Value = 131
$ws.Range("H6").Value = 684

# Irlanda's case counts overtook Japon's, so Irlanda moves above Japon
# in the sorted-by-total list (row 28 becomes Irlanda, row 29 becomes Japon).
# Row 28: now Irlanda, with its newly updated figures
$ws.Range("A28").Value = "Irlanda"
$ws.Range("B28").Value = 1329
$ws.Range("C28").Value = 204
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 1317
$ws.Range("F28").Value = 29
$ws.Range("H28").Value = 7

# Row 29: now Japon, carrying the figures Japon had before the update
$ws.Range("A29").Value = "Japon"
$ws.Range("B29").Value = 1193
$ws.Range("C29").Value = 65
$ws.Range("D29").Value = 285
$ws.Range("E29").Value = 865
$ws.Range("F29").Value = 54
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 43

# Row 34 - updated case counts
$ws.Range("B34").Value = 901
$ws.Range("C34").Value = 152
$ws.Range("E34").Value = 890

# Row 108 - updated case counts
$ws.Range("B108").Value = 57
$ws.Range("C108").Value = 6
$ws.Range("E108").Value = 57
